# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets carry the same table contents and need the identical set of
# updates, as per the commit's regenerated-output diff.

$wb = $excel.ActiveWorkbook

# Row number -> new F-column value (same update applies to both sheets).
$updates = @{
    3  = 98
    5  = 53
    7  = 60
    8  = 2039
    9  = 68
    10 = 104
    11 = 4424
    13 = 282
    14 = 103
    16 = 121
    19 = 77
    20 = 3267
    21 = 73
    22 = 489
    25 = 79
    26 = 91
    32 = 600
    33 = 1879
    34 = 292
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
